# Update the "想去人数" (column F) figures on the "展览" and "全部类型"
# worksheets to reflect the latest generated output.

$wb = $excel.ActiveWorkbook

# Row -> new value for column F
$updates = @{
    2  = 161
    3  = 7169
    4  = 5241
    9  = 103
    11 = 91
    12 = 194
    13 = 633
    14 = 215
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
